# Commit: "check if excel cell is empty"
# Populate the "vol min" / "vol max" columns (C/D) for the structures that
# now have computed values, and leave behind the block of blank,
# newly-formatted cells (columns G:O) that Excel's paste/insert left on
# several rows while the author was checking which cells were empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clinical Structures")

# --- vol min (C) / vol max (D) values -------------------------------------------------
$ws.Range("C10").Value = 112.02
$ws.Range("D10").Value = 5569.9416345124491

$ws.Range("C15").Value = 39.14
$ws.Range("D15").Value = 560.73982355024589

$ws.Range("C17").Value = 30.94
$ws.Range("D17").Value = 130.38317839653934

$ws.Range("C19").Value = 9.1
$ws.Range("D19").Value = 309.59451169159098

$ws.Range("C20").Value = 96.808165954476038
$ws.Range("D20").Value = 270.90074153704302

$ws.Range("C21").Value = 100.38
$ws.Range("D21").Value = 271.54249032790631

$ws.Range("C23").Value = 91.37
$ws.Range("D23").Value = 581.80274094704828

# --- leftover empty, but styled, cells (G:O) on the probed rows -----------------------
# These cells came along for the ride from a pasted block used while checking
# which structures still had an empty min/max cell; they carry no value, only
# a (visually default) cell format, which is why they still show up as real
# cells in the sheet's used range instead of being completely absent.
$ws.Range("G10:M12").NumberFormat = "General"
$ws.Range("G13:I14").NumberFormat = "General"
$ws.Range("G15:I19").NumberFormat = "General"
$ws.Range("I20:O20").NumberFormat = "General"
$ws.Range("I21:O22").NumberFormat = "General"

# --- view bookkeeping ------------------------------------------------------------------
$ws.Range("D23").Select()
$excel.ActiveWindow.Zoom = 100
